$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet1 rows 8-42 previously listed OIDs sequentially; four OIDs (rows 9, 10,
# 12 and 13) are pulled out and appended after row 38, with research notes
# added in column D, and everything else shifts up to fill the gap.

$ws.Range("A8").Value = "2.16.840.1.113883.1.11.11526"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = "https://terminology.hl7.org/3.1.0/ValueSet-v3-HumanLanguage.html"

$ws.Range("A9").Value = "2.16.840.1.113883.3.26.1.1"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = "https://terminology.hl7.org/2.0.0/CodeSystem-v3-nciThesaurus.html"

$ws.Range("A10").Value = "2.16.840.1.113883.3.88.12.3221.5.2"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0

$ws.Range("A11").Value = "2.16.840.1.113883.3.88.12.80.1"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0

$ws.Range("A12").Value = "2.16.840.1.113883.3.88.12.80.2"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0

$ws.Range("A13").Value = "2.16.840.1.113883.3.88.12.80.33"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0

$ws.Range("A14").Value = "2.16.840.1.113883.3.88.12.80.63"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0

$ws.Range("A15").Value = "2.16.840.1.113883.4.642.4.1131"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0

$ws.Range("A16").Value = "2.16.840.1.113883.5.1001"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0

$ws.Range("A17").Value = "2.16.840.1.113883.5.1002"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0

$ws.Range("A18").Value = "2.16.840.1.113883.5.1008"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0

$ws.Range("A19").Value = "2.16.840.1.113883.5.110"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0

$ws.Range("A20").Value = "2.16.840.1.113883.5.111"
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0

$ws.Range("A21").Value = "2.16.840.1.113883.5.14"
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0

$ws.Range("A22").Value = "2.16.840.1.113883.5.4"
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0

$ws.Range("A23").Value = "2.16.840.1.113883.5.41"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0

$ws.Range("A24").Value = "2.16.840.1.113883.5.6"
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 0

$ws.Range("A25").Value = "2.16.840.1.113883.5.88"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0

$ws.Range("A26").Value = "2.16.840.1.113883.5.89"
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0

$ws.Range("A27").Value = "2.16.840.1.113883.5.90"
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0

$ws.Range("A28").Value = "2.16.840.1.113883.6.1"
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = "https://terminology.hl7.org/CodeSystem-v3-loinc.html"

$ws.Range("A29").Value = "2.16.840.1.113883.6.101"
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = "https://terminology.hl7.org/CodeSystem-v3-nuccProviderCodes.html"

$ws.Range("A30").Value = "2.16.840.1.113883.6.104"
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0

$ws.Range("A31").Value = "2.16.840.1.113883.6.12"
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = "https://terminology.hl7.org/CodeSystem-CPT.html"

$ws.Range("A32").Value = "2.16.840.1.113883.6.13"
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = "https://terminology.hl7.org/CodeSystem-CDT.html"

$ws.Range("A33").Value = "2.16.840.1.113883.6.254"
$ws.Range("B33").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("D33").ClearContents()

$ws.Range("A34").Value = "2.16.840.1.113883.6.301.5"
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0

$ws.Range("A35").Value = "2.16.840.1.113883.6.4"
$ws.Range("B35").Value = 0
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = "https://terminology.hl7.org/CodeSystem-icd10PCS.html"

$ws.Range("A36").Value = "2.16.840.1.113883.6.90"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = "https://terminology.hl7.org/CodeSystem-icd10CM.html"

$ws.Range("A37").Value = "2.16.840.1.113883.6.96"
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = "https://terminology.hl7.org/CodeSystem-v3-snomed-CT.html"

$ws.Range("A38").Value = "2.16.840.1.113762.1.4.1247.YY"
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = "Invalid"

$ws.Range("A39:D39").ClearContents()

$ws.Range("A40").Value = "2.16.840.1.113883.11.20.9.28"
$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = "n/a"

$ws.Range("A41").Value = "2.16.840.1.113883.12.112"
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = "n/a"

$ws.Range("A42").Value = "2.16.840.1.113883.3.3719"
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = "Root for device IDs"

$ws.Range("A43").Value = "2.16.840.1.113883.3.5019.1.1"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 0
